$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Range("H16").Value = 15271.8
$ws_ALC.Range("I16").Value = 1675
$ws_ALC.Range("J16").Value = 24336.334
$ws_ALC.Range("K16").Value = 1675
$ws_ALC.Range("L16").Value = 24336.334
$ws_ALC.Range("M16").Value = -1445
$ws_ALC.Range("N16").Value = -24796.334

$ws_ALC.Range("H40").Value = 1200.3334
$ws_ALC.Range("I40").Value = 1100.5
$ws_ALC.Range("J40").Value = 1400
$ws_ALC.Range("K40").Value = 1100.5
$ws_ALC.Range("L40").Value = 1400
$ws_ALC.Range("M40").Value = -925.5
$ws_ALC.Range("N40").Value = -1750

$ws_ALC.Range("H41").Value = 749.2857
$ws_ALC.Range("I41").Value = 355.7143
$ws_ALC.Range("J41").Value = 1142.8572
$ws_ALC.Range("K41").Value = 355.7143
$ws_ALC.Range("L41").Value = 1142.8572
$ws_ALC.Range("M41").Value = 84.28570000000002
$ws_ALC.Range("N41").Value = -2022.8572

$ws_ALC.Range("H75").Value = 26599.666
$ws_ALC.Range("J75").Value = 26599.666
$ws_ALC.Range("L75").Value = 26599.666
$ws_ALC.Range("N75").Value = -28471.666

$ws_ALC.Range("H78").Value = 26599.666
$ws_ALC.Range("J78").Value = 26599.666
$ws_ALC.Range("L78").Value = 79798.99800000001
$ws_ALC.Range("N78").Value = -89158.99800000001

$ws_ALC.Range("H113").Value = 12917
$ws_ALC.Range("I113").Value = 3668.3333
$ws_ALC.Range("J113").Value = 15999.889
$ws_ALC.Range("K113").Value = 3668.3333
$ws_ALC.Range("L113").Value = 15999.889
$ws_ALC.Range("M113").Value = -414.3332999999998
$ws_ALC.Range("N113").Value = -22507.889

$ws_ALC.Range("H123").Value = 43113.332
$ws_ALC.Range("J123").Value = 43113.332
$ws_ALC.Range("L123").Value = 43113.332
$ws_ALC.Range("N123").Value = -52913.332

$ws_ALC.Range("H132").Value = 32264200
$ws_ALC.Range("I132").Value = 37043080
$ws_ALC.Range("J132").Value = 6750
$ws_ALC.Range("K132").Value = 111129240
$ws_ALC.Range("L132").Value = 20250
$ws_ALC.Range("M132").Value = -111126710
$ws_ALC.Range("N132").Value = -25310

$ws_ALC.Range("H138").Value = 3817.8289
$ws_ALC.Range("I138").Value = 1726.6875
$ws_ALC.Range("J138").Value = 4375.467
$ws_ALC.Range("K138").Value = 5180.0625
$ws_ALC.Range("L138").Value = 13126.401
$ws_ALC.Range("M138").Value = -40.0625
$ws_ALC.Range("N138").Value = -23406.401

$ws_ARM.Range("H7").Value = 39690
$ws_ARM.Range("J7").Value = 39690
$ws_ARM.Range("L7").Value = 39690
$ws_ARM.Range("N7").Value = -39918

$ws_ARM.Range("H32").Value = 11299.171
$ws_ARM.Range("I32").Value = 7493.3887
$ws_ARM.Range("K32").Value = 7493.3887
$ws_ARM.Range("M32").Value = -7206.3887

$ws_ARM.Range("H109").Value = 26500
$ws_ARM.Range("J109").Value = 26500
$ws_ARM.Range("L109").Value = 26500
$ws_ARM.Range("N109").Value = -29274

$ws_ARM.Range("H122").Value = 3049.9443
$ws_ARM.Range("I122").Value = 1658.25
$ws_ARM.Range("K122").Value = 4974.75
$ws_ARM.Range("M122").Value = -2524.75

$ws_BSM.Range("H20").Value = 3537.0557
$ws_BSM.Range("I20").Value = 4060.6365
$ws_BSM.Range("J20").Value = 2714.2856
$ws_BSM.Range("K20").Value = 4060.6365
$ws_BSM.Range("L20").Value = 2714.2856
$ws_BSM.Range("M20").Value = -3813.6365
$ws_BSM.Range("N20").Value = -3208.2856

$ws_CRP.Range("H31").Value = 4818
$ws_CRP.Range("I31").Value = 2318
$ws_CRP.Range("K31").Value = 2318
$ws_CRP.Range("M31").Value = -2023

$ws_CRP.Range("H33").Value = 5983
$ws_CRP.Range("I33").Value = 5983
$ws_CRP.Range("J33").Value = 0
$ws_CRP.Range("K33").Value = 5983
$ws_CRP.Range("L33").Value = 0
$ws_CRP.Range("M33").Value = -5604
$ws_CRP.Range("N33").ClearContents()

$ws_CRP.Range("H34").Value = 4818
$ws_CRP.Range("I34").Value = 2318
$ws_CRP.Range("K34").Value = 2318
$ws_CRP.Range("M34").Value = -2116

$ws_CRP.Range("H39").Value = 17064.438
$ws_CRP.Range("I39").Value = 2974.8333
$ws_CRP.Range("J39").Value = 25518.2
$ws_CRP.Range("K39").Value = 2974.8333
$ws_CRP.Range("L39").Value = 25518.2
$ws_CRP.Range("M39").Value = -2583.8333
$ws_CRP.Range("N39").Value = -26300.2

$ws_CRP.Range("H49").Value = 17064.438
$ws_CRP.Range("I49").Value = 2974.8333
$ws_CRP.Range("J49").Value = 25518.2
$ws_CRP.Range("K49").Value = 2974.8333
$ws_CRP.Range("L49").Value = 25518.2
$ws_CRP.Range("M49").Value = -2792.8333
$ws_CRP.Range("N49").Value = -25882.2

$ws_CRP.Range("H59").Value = 28844.75
$ws_CRP.Range("J59").Value = 28844.75
$ws_CRP.Range("L59").Value = 28844.75
$ws_CRP.Range("N59").Value = -31134.75

$ws_CRP.Range("H60").Value = 27157
$ws_CRP.Range("J60").Value = 29174
$ws_CRP.Range("L60").Value = 29174
$ws_CRP.Range("N60").Value = -30196

$ws_CRP.Range("H82").Value = 43000
$ws_CRP.Range("J82").Value = 43000
$ws_CRP.Range("L82").Value = 43000
$ws_CRP.Range("N82").Value = -43722

$ws_CRP.Range("H85").Value = 43000
$ws_CRP.Range("J85").Value = 43000
$ws_CRP.Range("L85").Value = 43000
$ws_CRP.Range("N85").Value = -45496

$ws_CRP.Range("H99").Value = 3981.4285
$ws_CRP.Range("I99").Value = 1888.75
$ws_CRP.Range("J99").Value = 6771.6665
$ws_CRP.Range("K99").Value = 1888.75
$ws_CRP.Range("L99").Value = 6771.6665
$ws_CRP.Range("M99").Value = -390.75
$ws_CRP.Range("N99").Value = -9767.666499999999

$ws_CRP.Range("H126").Value = 3981.4285
$ws_CRP.Range("I126").Value = 1888.75
$ws_CRP.Range("J126").Value = 6771.6665
$ws_CRP.Range("K126").Value = 5666.25
$ws_CRP.Range("L126").Value = 20314.9995
$ws_CRP.Range("M126").Value = -3196.25
$ws_CRP.Range("N126").Value = -25254.9995

$ws_CRP.Range("H139").Value = 43004.285
$ws_CRP.Range("J139").Value = 43004.285
$ws_CRP.Range("L139").Value = 43004.285
$ws_CRP.Range("N139").Value = -53284.285

$ws_CUL.Range("H5").Value = 1638.8
$ws_CUL.Range("I5").Value = 415.61905
$ws_CUL.Range("J5").Value = 3473.5715
$ws_CUL.Range("K5").Value = 1246.85715
$ws_CUL.Range("L5").Value = 10420.7145
$ws_CUL.Range("M5").Value = -1134.85715
$ws_CUL.Range("N5").Value = -10644.7145

$ws_CUL.Range("H25").Value = 4000.25
$ws_CUL.Range("I25").Value = 1001
$ws_CUL.Range("K25").Value = 3003
$ws_CUL.Range("M25").Value = -2834

$ws_CUL.Range("H30").Value = 4000.25
$ws_CUL.Range("I30").Value = 1001
$ws_CUL.Range("K30").Value = 3003
$ws_CUL.Range("M30").Value = -2901

$ws_CUL.Range("H113").Value = 556.3182
$ws_CUL.Range("I113").Value = 611.17645
$ws_CUL.Range("J113").Value = 521.7778
$ws_CUL.Range("K113").Value = 1833.52935
$ws_CUL.Range("L113").Value = 1565.3334
$ws_CUL.Range("M113").Value = 336.4706499999998
$ws_CUL.Range("N113").Value = -5905.3334

$ws_CUL.Range("H122").Value = 3264.9783
$ws_CUL.Range("J122").Value = 3618.5
$ws_CUL.Range("L122").Value = 32566.5
$ws_CUL.Range("N122").Value = -37466.5

$ws_CUL.Range("H123").Value = 2342.7144
$ws_CUL.Range("J123").Value = 2149.75
$ws_CUL.Range("L123").Value = 6449.25
$ws_CUL.Range("N123").Value = -11349.25

$ws_CUL.Range("H135").Value = 1638.8
$ws_CUL.Range("I135").Value = 415.61905
$ws_CUL.Range("J135").Value = 3473.5715
$ws_CUL.Range("K135").Value = 3740.57145
$ws_CUL.Range("L135").Value = 31262.1435
$ws_CUL.Range("M135").Value = -1205.57145
$ws_CUL.Range("N135").Value = -36332.1435

$ws_CUL.Range("H140").Value = 2048.7083
$ws_CUL.Range("I140").Value = 1573.0625
$ws_CUL.Range("K140").Value = 4719.1875
$ws_CUL.Range("M140").Value = 460.8125

$ws_GSM.Range("H70").Value = 6256.926
$ws_GSM.Range("I70").Value = 5580.0557
$ws_GSM.Range("K70").Value = 5580.0557
$ws_GSM.Range("M70").Value = -5310.0557

$ws_GSM.Range("H73").Value = 6256.926
$ws_GSM.Range("I73").Value = 5580.0557
$ws_GSM.Range("K73").Value = 5580.0557
$ws_GSM.Range("M73").Value = -4644.0557

$ws_GSM.Range("H102").Value = 2142.7678
$ws_GSM.Range("I102").Value = 1819.0834
$ws_GSM.Range("J102").Value = 4084.875
$ws_GSM.Range("K102").Value = 1819.0834
$ws_GSM.Range("L102").Value = 4084.875
$ws_GSM.Range("M102").Value = -197.0834
$ws_GSM.Range("N102").Value = -7328.875

$ws_GSM.Range("H122").Value = 2800.9092
$ws_GSM.Range("I122").Value = 2158.0527
$ws_GSM.Range("J122").Value = 3673.3572
$ws_GSM.Range("K122").Value = 6474.158100000001
$ws_GSM.Range("L122").Value = 11020.0716
$ws_GSM.Range("M122").Value = -4024.158100000001
$ws_GSM.Range("N122").Value = -15920.0716

$ws_GSM.Range("H132").Value = 2804.5625
$ws_GSM.Range("I132").Value = 630.53845
$ws_GSM.Range("J132").Value = 4292.0527
$ws_GSM.Range("K132").Value = 1891.61535
$ws_GSM.Range("L132").Value = 12876.1581
$ws_GSM.Range("M132").Value = 638.38465
$ws_GSM.Range("N132").Value = -17936.1581

$ws_LTW.Range("H22").Value = 11907095
$ws_LTW.Range("I22").Value = 22729210
$ws_LTW.Range("J22").Value = 2770.3
$ws_LTW.Range("K22").Value = 22729210
$ws_LTW.Range("L22").Value = 2770.3
$ws_LTW.Range("M22").Value = -22728915
$ws_LTW.Range("N22").Value = -3360.3

$ws_LTW.Range("H27").Value = 11907095
$ws_LTW.Range("I27").Value = 22729210
$ws_LTW.Range("J27").Value = 2770.3
$ws_LTW.Range("K27").Value = 22729210
$ws_LTW.Range("L27").Value = 2770.3
$ws_LTW.Range("M27").Value = -22729103
$ws_LTW.Range("N27").Value = -2984.3

$ws_LTW.Range("H55").Value = 339.21738
$ws_LTW.Range("I55").Value = 284
$ws_LTW.Range("J55").Value = 399.45456
$ws_LTW.Range("K55").Value = 284
$ws_LTW.Range("L55").Value = 399.45456
$ws_LTW.Range("M55").Value = -111
$ws_LTW.Range("N55").Value = -745.45456

$ws_LTW.Range("H132").Value = 5197.385
$ws_LTW.Range("I132").Value = 1255.6471
$ws_LTW.Range("J132").Value = 12642.889
$ws_LTW.Range("K132").Value = 3766.9413
$ws_LTW.Range("L132").Value = 37928.667
$ws_LTW.Range("M132").Value = -1236.9413
$ws_LTW.Range("N132").Value = -42988.667

$ws_LTW.Range("H136").Value = 3105.9678
$ws_LTW.Range("I136").Value = 1615.8334
$ws_LTW.Range("K136").Value = 4847.5002
$ws_LTW.Range("M136").Value = -2297.5002

$ws_WVR.Range("H136").Value = 2957.4546
$ws_WVR.Range("I136").Value = 703.5454999999999
$ws_WVR.Range("J136").Value = 7465.273
$ws_WVR.Range("K136").Value = 2110.6365
$ws_WVR.Range("L136").Value = 22395.819
$ws_WVR.Range("M136").Value = 439.3635000000004
$ws_WVR.Range("N136").Value = -27495.819
